# Added new spider monsters and new 2D map 'abandoned hut'
#
# - Adds two new monster/monster-group entries to the "Monsters" sheet:
#     row 3: ID 58 "Höhlenspinne" in "Monster in Ship's end",
#            Group ID 88 "3x Höhlenspinne"
#     row 4: Group ID 89 "4x Höhlenspinne"
# - Makes "Monsters" the active sheet/tab (previously "Todo" was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monsters")

$ws.Cells.Item(3, 1).Value = 58
$ws.Cells.Item(3, 2).Value = "Höhlenspinne"
$ws.Cells.Item(3, 3).Value = "Monster in Ship's end"
$ws.Cells.Item(3, 8).Value = 88
$ws.Cells.Item(3, 9).Value = "3x Höhlenspinne"

$ws.Cells.Item(4, 8).Value = 89
$ws.Cells.Item(4, 9).Value = "4x Höhlenspinne"

[void]$ws.Range("J4").Select()
$ws.Activate()
